# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The employee/period detail table (rows 16-48, columns C:F) is re-sorted:
# originally grouped by worker (each worker's periods listed together),
# now grouped by period (each period's workers listed together), while
# keeping the exact same set of (worker, period) -> Valor Mora pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora, F = Valor Mora
$rows = @(
    @{ Row = 16; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2308"; Valor = 46400 },
    @{ Row = 17; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2309"; Valor = 46400 },
    @{ Row = 18; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2310"; Valor = 46400 },
    @{ Row = 19; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2311"; Valor = 46400 },
    @{ Row = 20; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2311"; Valor = 46400 },
    @{ Row = 21; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2312"; Valor = 46400 },
    @{ Row = 22; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2312"; Valor = 46400 },
    @{ Row = 23; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2312"; Valor = 46400 },
    @{ Row = 24; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2312"; Valor = 46400 },
    @{ Row = 25; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2401"; Valor = 46400 },
    @{ Row = 26; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2401"; Valor = 46400 },
    @{ Row = 27; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2401"; Valor = 46400 },
    @{ Row = 28; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2401"; Valor = 46400 },
    @{ Row = 29; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2402"; Valor = 46400 },
    @{ Row = 30; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2402"; Valor = 46400 },
    @{ Row = 31; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2402"; Valor = 46400 },
    @{ Row = 32; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2402"; Valor = 46400 },
    @{ Row = 33; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2403"; Valor = 46400 },
    @{ Row = 34; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2403"; Valor = 46400 },
    @{ Row = 35; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2403"; Valor = 46400 },
    @{ Row = 36; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2403"; Valor = 46400 },
    @{ Row = 37; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2404"; Valor = 46400 },
    @{ Row = 38; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2404"; Valor = 46400 },
    @{ Row = 39; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2404"; Valor = 46400 },
    @{ Row = 40; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2404"; Valor = 46400 },
    @{ Row = 41; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2405"; Valor = 46400 },
    @{ Row = 42; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2405"; Valor = 46400 },
    @{ Row = 43; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2405"; Valor = 46400 },
    @{ Row = 44; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2405"; Valor = 46400 },
    @{ Row = 45; Doc = "1049828873"; Nombre = "LUIS GUILLERMO CANTILLO FLOREZ";     Periodo = "2406"; Valor = 18560 },
    @{ Row = 46; Doc = "1143329466"; Nombre = "HERIBERTO DE JESUS MARQUEZ SOLIS";   Periodo = "2406"; Valor = 18560 },
    @{ Row = 47; Doc = "23139700";   Nombre = "YESEIRA SUAREZ BATISTA";             Periodo = "2406"; Valor = 18560 },
    @{ Row = 48; Doc = "1051889025"; Nombre = "MARILUZ CONEO JIMENEZ";              Periodo = "2406"; Valor = 18560 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Valor
}
